# Daily automatic data refresh for the EPEX Spot prices workbook.
# Adds the "15-aug" column (BK) to the "Prix Spot" sheet, and appends the
# 2025-08-13 row to the "Gaz" and "CO2" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Prix Spot" ---------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Header cell BK1: copy BJ1 (keeps the bold/centered/bordered header
# style) onto BK1, then overwrite its text.
$wsSpot.Range("BJ1").Copy($wsSpot.Range("BK1"))
$wsSpot.Range("BK1").Value = "15-aug"

# Hourly price values for the new day, row by row.
$spotValues = @{
    2  = 90.52
    3  = 83.53
    4  = 79.41
    5  = 68.47
    6  = 76.09
    7  = 65.65000000000001
    8  = 73.45999999999999
    9  = 79.98999999999999
    10 = 71.89
    11 = 43.97
    12 = 25
    13 = 3
    14 = 0.08
    15 = 0
    16 = 0
    17 = 1.72
    18 = 23.53
    19 = 33.23
    20 = 61.12
    21 = 100.82
    22 = 111.39
    23 = 102
    24 = 102.89
    25 = 93.45
}

foreach ($row in $spotValues.Keys) {
    $wsSpot.Range("BK$row").Value = $spotValues[$row]
}

# --- Sheet 2: "Gaz" ------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A60").NumberFormat = "@"
$wsGaz.Range("A60").Value = "2025-08-13"
$wsGaz.Range("A60").Style = "Normal"
$wsGaz.Range("B60").Value = 31.85

# --- Sheet 3: "CO2" -------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A60").NumberFormat = "@"
$wsCo2.Range("A60").Value = "2025-08-13"
$wsCo2.Range("A60").Style = "Normal"
$wsCo2.Range("B60").Value = 71.06999999999999
